# Weekly refresh of the "Fruta / hortaliza" consolidated sheet: the source
# feed re-shuffled which data row each market record landed on. Columns
# D (Fecha), K (Variedad), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), Q (Unidad de
# comercializacion), R (Origen) and S (Precio $/Kg) move between rows
# 2-29; every other column is identical across rows so it is untouched.
#
# mapping[i] = the 1-based source row (in the ORIGINAL sheet) whose
# D/K/L/M/N/O/P/Q/R/S values now belong on destination row (i + 2).
$mapping = @(29, 13, 19, 25, 27, 20, 12, 3, 28, 26, 5, 17, 7, 6, 8, 9, 2, 21, 10, 18, 11, 23, 4, 14, 15, 16, 22, 24)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 29
$cols = @("D", "K", "L", "M", "N", "O", "P", "Q", "R", "S")

# 1) Snapshot the current (pre-edit) values for every touched column/row
#    before any writes happen, so the shuffle reads only original data.
#    (Value2 is used instead of Value - chaining .Value straight into
#    another property assignment does not resolve through this host.)
$snapshot = @{}
foreach ($col in $cols) {
    $colValues = @{}
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $colValues[$r] = $ws.Range("$col$r").Value2
    }
    $snapshot[$col] = $colValues
}

# 2) Write the snapshotted values back out according to the mapping.
for ($i = 0; $i -lt $mapping.Length; $i++) {
    $destRow = $firstRow + $i
    $srcRow = $mapping[$i]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $snapshot[$col][$srcRow]
    }
}
